# Auto-generated edit script applying the cryptos.xlsx symbol-list update
# (commit: "Updated symbol list on Tue Jan 31 22:48:55 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '311.37'
Set-TextCell 2 5 '1.84%'
Set-TextCell 3 5 '0.56%'
Set-TextCell 4 4 '5.122'
Set-TextCell 4 5 '0.73%'
Set-TextCell 5 4 '0.07836'
Set-TextCell 5 5 '1.61%'
Set-TextCell 6 2 'GateToken'
Set-TextCell 6 3 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell 6 4 '4.429'
Set-TextCell 6 5 '1.86%'
Set-TextCell 7 2 'FTXToken'
Set-TextCell 7 3 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell 7 4 '1.914'
Set-TextCell 7 5 '1.94%'
Set-TextCell 8 2 'KuCoinToken'
Set-TextCell 8 3 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextCell 8 4 '8.254'
Set-TextCell 8 5 '1.03%'
Set-TextCell 9 2 'BTSEToken'
Set-TextCell 9 3 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextCell 9 4 '2.952'
Set-TextCell 9 5 '-6.78%'
Set-TextCell 10 2 'MXToken'
Set-TextCell 10 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 10 4 '0.9193'
Set-TextCell 11 2 'LiechtensteinCryptoassetsExchange'
Set-TextCell 11 3 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell 11 4 '0.1197'
Set-TextCell 11 5 '-0.20%'
Set-TextCell 12 2 'WazirX'
Set-TextCell 12 3 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell 12 4 '0.1914'
Set-TextCell 12 5 '2.73%'
Set-TextCell 13 2 'MandalaExchangeToken'
Set-TextCell 13 3 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell 13 4 '0.08986'
Set-TextCell 13 5 '2.65%'
Set-TextCell 14 2 'BitrueCoin'
Set-TextCell 14 3 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell 14 4 '0.03351'
Set-TextCell 14 5 '-1.03%'
Set-TextCell 15 2 'BitMartToken'
Set-TextCell 15 3 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell 15 4 '0.09594'
Set-TextCell 15 5 '-1.00%'
Set-TextCell 16 2 'BitForexToken'
Set-TextCell 16 3 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell 16 4 '0.001377'
Set-TextCell 16 5 '0.86%'
Set-TextCell 17 2 'TigerCash'
Set-TextCell 17 3 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell 17 4 '0.005725'
Set-TextCell 17 5 '-3.69%'
Set-TextCell 18 2 'LEO'
Set-TextCell 18 3 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell 18 4 '3.538'
Set-TextCell 18 5 '-1.81%'
Set-TextCell 19 4 '0.3441'
Set-TextCell 19 5 '0.95%'
Set-TextCell 20 4 '5.247'
Set-TextCell 20 5 '4.63%'
Set-TextCell 21 4 '0.1285'
Set-TextCell 21 5 '0.62%'
Set-TextCell 22 4 '0.2593'
Set-TextCell 22 5 '-0.11%'
Set-TextCell 23 4 '0.04356'
Set-TextCell 23 5 '0.74%'
Set-TextCell 24 5 '3.06%'
Set-TextCell 25 4 '0.004657'
Set-TextCell 25 5 '10.43%'
Set-TextCell 26 4 '0.0001360'
Set-TextCell 26 5 '0.57%'
Set-TextCell 27 4 '0.0003996'
Set-TextCell 27 5 '-98.10%'
Set-TextCell 39 4 '0.02261'
Set-TextCell 39 5 '3.86%'
Set-TextCell 40 4 '0.05052'
Set-TextCell 40 5 '3.45%'
Set-TextCell 41 4 '0.007462'
Set-TextCell 41 5 '-1.39%'
Set-TextCell 42 4 '0.009055'
Set-TextCell 42 5 '-8.84%'
Set-TextCell 43 4 '0.1348'
Set-TextCell 43 5 '0.84%'
Set-TextCell 44 4 '0.001950'
Set-TextCell 44 5 '-2.29%'
Set-TextCell 45 4 '0.009303'
Set-TextCell 45 5 '2.10%'
Set-TextCell 46 4 '0.00006565'
Set-TextCell 46 5 '0.24%'
Set-TextCell 47 5 '-0.16%'
Set-TextCell 48 2 'CoinbaseStockToken'
Set-TextCell 48 3 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextCell 48 4 '0.001001'
Set-TextCell 48 5 '-23.08%'
Set-TextCell 49 2 'BOLO'
Set-TextCell 49 3 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextCell 49 4 '0.003369'
Set-TextCell 49 5 '12.18%'
Set-TextCell 50 4 '0.00002100'
Set-TextCell 50 5 '-0.16%'
Set-TextCell 51 4 '0.0002000'
Set-TextCell 51 5 '-0.16%'
